# Informe-02-020005-A-TC: rework the DSD header block.
#  - Row 1: human-readable (capitalized) column labels instead of raw codelist ids
#  - Row 2: iaest-measure: URNs, reordered to align with the new row-1 labels
#  - Row 3: "medida" labels (unchanged)
#  - Row 4: data types per measure - all xsd:string except "Municipio" (column C),
#           which is xsd:int (fixes the erroneous measure generation)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the cells we are about to rewrite so stale shared-string references
# don't linger.
$ws.Range("A1:G2").ClearContents()
$ws.Range("A4:G4").ClearContents()

# Row 1 - column headers
$ws.Range("A1").Value = "Municipio"
$ws.Range("B1").Value = "Entidad singular"
$ws.Range("C1").Value = "Personas"
$ws.Range("D1").Value = "Núcleo"
$ws.Range("E1").Value = "Comarca"
$ws.Range("F1").Value = "Provincia"
$ws.Range("G1").Value = "Núcleo/diseminado"

# Row 2 - iaest-measure URNs matching the row-1 headers
$ws.Range("A2").Value = "iaest-measure:municipio"
$ws.Range("B2").Value = "iaest-measure:entidad-singular"
$ws.Range("C2").Value = "iaest-measure:personas"
$ws.Range("D2").Value = "iaest-measure:nucleo"
$ws.Range("E2").Value = "iaest-measure:comarca"
$ws.Range("F2").Value = "iaest-measure:provincia"
$ws.Range("G2").Value = "iaest-measure:nucleodiseminado"

# Row 3 ("medida") is untouched.

# Row 4 - data types: everything is xsd:string except Municipio (xsd:int)
$ws.Range("A4").Value = "xsd:string"
$ws.Range("B4").Value = "xsd:string"
$ws.Range("C4").Value = "xsd:int"
$ws.Range("D4").Value = "xsd:string"
$ws.Range("E4").Value = "xsd:string"
$ws.Range("F4").Value = "xsd:string"
$ws.Range("G4").Value = "xsd:string"
